$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# End Year: 2030 -> 2050
$ws.Range("B4").Value = 2050

# maximum_investment_capacity_per_year: 300 -> 1000000
$ws.Range("B13").Value = 1000000

# realistic_candidate_capacities_tobe_installed: TRUE -> FALSE
$ws.Range("B16").Value = $false

# realistic_candidate_capacities_for_future: TRUE -> FALSE
$ws.Range("B17").Value = $false

# dummy_capacity: 100 -> 300
$ws.Range("B19").Value = 300
